$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column F: header "11_03_2024" and the new sales/agendamiento figures
$ws.Range("F1").Value = "11_03_2024"
$ws.Range("F2").Value = 1070
$ws.Range("F3").Value = 1037
$ws.Range("F4").Value = 1548
$ws.Range("F5").Value = 2930
$ws.Range("F6").Value = 141

$ws.Range("F6").Select()
